$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25
$ws.Cells.Item($row, 1).Value = "Cudahy"
$ws.Cells.Item($row, 2).Value = "Chuds"
$ws.Cells.Item($row, 3).Value = "The Resistance"
$ws.Cells.Item($row, 4).Value = "penguino#2114"
$ws.Cells.Item($row, 5).Value = "Pingu"
$ws.Cells.Item($row, 6).Value = "air raid"
$ws.Cells.Item($row, 7).Value = "5-2"
$ws.Cells.Item($row, 8).Value = "0-0"
